$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Number of existing units" (column F) for rows with a new value
$ws.Range("F8").Value = 19
$ws.Range("F10").Value = 2
$ws.Range("F12").Value = 23
$ws.Range("F14").Value = 10
$ws.Range("F16").Value = 83

# Update "Maximum number of units that can be invested in" (column I) for rows 8-18
$ws.Range("I8").Value = 8
$ws.Range("I9").Value = 8
$ws.Range("I10").Value = 8
$ws.Range("I11").Value = 8
$ws.Range("I12").Value = 8
$ws.Range("I13").Value = 8
$ws.Range("I14").Value = 8
$ws.Range("I15").Value = 8
$ws.Range("I16").Value = 8
$ws.Range("I17").Value = 8
$ws.Range("I18").Value = 8

# Update the active cell / selection to H22
$ws.Range("H22").Select()
